# Build site at 2022-01-09 00:29:46 UTC
# Insert a "Docentes responsaveis:" block (5 rows) after row 11 ("Objectives:")
# by shifting the existing content of rows 12-22 down to rows 17-27, then
# filling the freed rows 12-16 with the new teacher list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteValues = -4163, xlPasteFormats = -4122
$xlPasteValues = -4163
$xlPasteFormats = -4122

# Columns that actually contain data for each source row (12-22), and the
# row height (in points) that each source row currently has (None -> no
# explicit row height).
$colsFor = @{
  12 = "A:C"; 13 = "A:A"; 14 = "A:C"; 15 = "A:A"; 16 = "A:A";
  17 = "A:C"; 18 = "A:C"; 19 = "A:C"; 20 = "A:C"; 21 = "A:A"; 22 = "B:C"
}

# Move rows 22 down to 12 (descending) so every source row is read before
# it is overwritten as somebody else's destination row.
for ($src = 22; $src -ge 12; $src--) {
    $dst = $src + 5
    $colRange = $colsFor[$src]
    $parts = $colRange.Split(":")
    $srcRange = $ws.Range($parts[0] + $src + ":" + $parts[1] + $src)
    $dstRange = $ws.Range($parts[0] + $dst + ":" + $parts[1] + $dst)

    $srcRange.Copy()
    $dstRange.PasteSpecial($xlPasteValues)
    $srcRange.Copy()
    $dstRange.PasteSpecial($xlPasteFormats)

    $srcHeight = $ws.Rows($src).RowHeight
    $ws.Rows($dst).RowHeight = $srcHeight
}
$ws.Application.CutCopyMode = $false

# Now clear out the old content left behind in rows 12-16 (it has already
# been copied down to rows 17-21) so we can place the new block there.
$ws.Range("A12:C16").ClearContents()
$ws.Rows("12:16").RowHeight = $ws.StandardHeight

# Row 12: "Docentes responsaveis:" (column A only, bold header style)
$ws.Range("A11").Copy()
$ws.Range("A12").PasteSpecial($xlPasteFormats)
$ws.Range("A12").Value = "Docentes responsáveis:"

# Rows 13-16: the four teacher names, in columns B and C only
$ws.Range("B17").Copy()
$ws.Range("B13:B16").PasteSpecial($xlPasteFormats)
$ws.Range("C17").Copy()
$ws.Range("C13:C16").PasteSpecial($xlPasteFormats)
$ws.Application.CutCopyMode = $false

$ws.Range("B13").Value = "471420 - Carlos Antonio Reis Pereira Baptista"
$ws.Range("C13").Value = "471420 - Carlos Antonio Reis Pereira Baptista"

$ws.Range("B14").Value = "3480026 - João Paulo Pascon"
$ws.Range("C14").Value = "3480026 - João Paulo Pascon"

$ws.Range("B15").Value = "5840793 - Sérgio Schneider"
$ws.Range("C15").Value = "5840793 - Sérgio Schneider"

$ws.Range("B16").Value = "7797767 - Viktor Pastoukhov"
$ws.Range("C16").Value = "7797767 - Viktor Pastoukhov"

Write-Host "done"
